$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", matching the formatting of the other
# header cells (B1:G1, style index 1 -> bold, bordered, centered/top).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data cell H2 = 0 (numeric)
$ws.Range("H2").Value = 0
